$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 220.76923
$ws.Range("I2").Value = 275.66666
$ws.Range("J2").Value = 97.25
$ws.Range("K2").Value = 275.66666
$ws.Range("L2").Value = 97.25
$ws.Range("M2").Value = -162.66666
$ws.Range("N2").Value = -323.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 794.6429000000001
$ws.Range("I19").Value = 703.1818
$ws.Range("J19").Value = 1130
$ws.Range("K19").Value = 703.1818
$ws.Range("L19").Value = 1130
$ws.Range("M19").Value = -528.1818
$ws.Range("N19").Value = -1480

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9999.933999999999
$ws.Range("I21").Value = 9999.799999999999
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 9999.799999999999
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -9531.799999999999
$ws.Range("N21").Value = -10936

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 9999.933999999999
$ws.Range("I23").Value = 9999.799999999999
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 9999.799999999999
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -9765.799999999999
$ws.Range("N23").Value = -10468

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 470.27274
$ws.Range("I38").Value = 130.33333
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 390.99999
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -18.99998999999997
$ws.Range("N38").Value = -6744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 426
$ws.Range("I58").Value = 160.83333
$ws.Range("K58").Value = 482.49999
$ws.Range("M58").Value = -332.49999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1721.2963
$ws.Range("I98").Value = 1830
$ws.Range("K98").Value = 1830
$ws.Range("M98").Value = -332

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1538.8043
$ws.Range("J112").Value = 1640.3572
$ws.Range("L112").Value = 4921.071599999999
$ws.Range("N112").Value = -7137.071599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 644.25
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7878.3335
$ws.Range("I116").Value = 1521
$ws.Range("J116").Value = 15825
$ws.Range("K116").Value = 1521
$ws.Range("L116").Value = 15825
$ws.Range("M116").Value = 1921
$ws.Range("N116").Value = -22709

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1721.2963
$ws.Range("I122").Value = 1830
$ws.Range("K122").Value = 5490
$ws.Range("M122").Value = -3040

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 37875
$ws.Range("J123").Value = 37875
$ws.Range("L123").Value = 37875
$ws.Range("N123").Value = -47675

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1326863.6
$ws.Range("I132").Value = 2673.5938
$ws.Range("J132").Value = 9801680
$ws.Range("K132").Value = 8020.7814
$ws.Range("L132").Value = 29405040
$ws.Range("M132").Value = -5490.7814
$ws.Range("N132").Value = -29410100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3113.375
$ws.Range("I141").Value = 2195.842
$ws.Range("K141").Value = 6587.526
$ws.Range("M141").Value = -1407.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 100000000
$ws.Range("I10").Value = 100000000
$ws.Range("K10").Value = 100000000
$ws.Range("M10").Value = -99999830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4347173
$ws.Range("I32").Value = 5146081
$ws.Range("J32").Value = 19757
$ws.Range("K32").Value = 5146081
$ws.Range("L32").Value = 19757
$ws.Range("M32").Value = -5145794
$ws.Range("N32").Value = -20331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 52738224
$ws.Range("I61").Value = 66734484
$ws.Range("J61").Value = 252250
$ws.Range("K61").Value = 66734484
$ws.Range("L61").Value = 252250
$ws.Range("M61").Value = -66734272
$ws.Range("N61").Value = -252674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2415.5
$ws.Range("I63").Value = 2415.5
$ws.Range("K63").Value = 2415.5
$ws.Range("M63").Value = -1729.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2415.5
$ws.Range("I66").Value = 2415.5
$ws.Range("K66").Value = 12077.5
$ws.Range("M66").Value = -8645.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4631813
$ws.Range("I122").Value = 2288.611
$ws.Range("J122").Value = 18520386
$ws.Range("K122").Value = 6865.833
$ws.Range("L122").Value = 55561158
$ws.Range("M122").Value = -4415.833
$ws.Range("N122").Value = -55566058

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 44307
$ws.Range("I132").Value = 31844.395
$ws.Range("J132").Value = 73683.14
$ws.Range("K132").Value = 95533.185
$ws.Range("L132").Value = 221049.42
$ws.Range("M132").Value = -93003.185
$ws.Range("N132").Value = -226109.42

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 52738224
$ws.Range("I136").Value = 66734484
$ws.Range("J136").Value = 252250
$ws.Range("K136").Value = 200203452
$ws.Range("L136").Value = 756750
$ws.Range("M136").Value = -200200902
$ws.Range("N136").Value = -761850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 30576.75
$ws.Range("J86").Value = 4753.5
$ws.Range("L86").Value = 4753.5
$ws.Range("N86").Value = -6999.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 30576.75
$ws.Range("J89").Value = 4753.5
$ws.Range("L89").Value = 23767.5
$ws.Range("N89").Value = -34999.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1097.9375
$ws.Range("I99").Value = 1071.1333
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1071.1333
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 426.8667
$ws.Range("N99").Value = -4496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11575.304
$ws.Range("I31").Value = 31689.834
$ws.Range("J31").Value = 2047.3684
$ws.Range("K31").Value = 31689.834
$ws.Range("L31").Value = 2047.3684
$ws.Range("M31").Value = -31394.834
$ws.Range("N31").Value = -2637.3684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11575.304
$ws.Range("I34").Value = 31689.834
$ws.Range("J34").Value = 2047.3684
$ws.Range("K34").Value = 31689.834
$ws.Range("L34").Value = 2047.3684
$ws.Range("M34").Value = -31487.834
$ws.Range("N34").Value = -2451.3684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3241.4
$ws.Range("I99").Value = 3601.6667
$ws.Range("J99").Value = 1350
$ws.Range("K99").Value = 3601.6667
$ws.Range("L99").Value = 1350
$ws.Range("M99").Value = -2103.6667
$ws.Range("N99").Value = -4346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3241.4
$ws.Range("I126").Value = 3601.6667
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 10805.0001
$ws.Range("L126").Value = 4050
$ws.Range("M126").Value = -8335.000100000001
$ws.Range("N126").Value = -8990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1929
$ws.Range("J97").Value = 861.2
$ws.Range("L97").Value = 2583.6
$ws.Range("N97").Value = -3575.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 854.8939
$ws.Range("I107").Value = 407.2449
$ws.Range("J107").Value = 2145.1765
$ws.Range("K107").Value = 1221.7347
$ws.Range("L107").Value = 6435.529500000001
$ws.Range("M107").Value = 698.2653
$ws.Range("N107").Value = -10275.5295

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3220.889
$ws.Range("I109").Value = 788
$ws.Range("J109").Value = 3525
$ws.Range("K109").Value = 2364
$ws.Range("L109").Value = 10575
$ws.Range("M109").Value = -1324
$ws.Range("N109").Value = -12655

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 5566.25
$ws.Range("I120").Value = 5566.25
$ws.Range("K120").Value = 16698.75
$ws.Range("M120").Value = -11860.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 813.7692
$ws.Range("J131").Value = 916.9
$ws.Range("L131").Value = 2750.7
$ws.Range("N131").Value = -12830.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12432.777
$ws.Range("J5").Value = 11486.875
$ws.Range("L5").Value = 11486.875
$ws.Range("N5").Value = -11710.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2178.889
$ws.Range("I97").Value = 2230
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 2230
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1734
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 49163
$ws.Range("J125").Value = 49163
$ws.Range("L125").Value = 49163
$ws.Range("N125").Value = -54083

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 34839.4
$ws.Range("J133").Value = 34839.4
$ws.Range("L133").Value = 34839.4
$ws.Range("N133").Value = -39899.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 56699.453
$ws.Range("I136").Value = 28846.379
$ws.Range("J136").Value = 121109.69
$ws.Range("K136").Value = 86539.137
$ws.Range("L136").Value = 363329.07
$ws.Range("M136").Value = -83989.137
$ws.Range("N136").Value = -368429.07

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2537501.5
$ws.Range("I2").Value = 5005000
$ws.Range("K2").Value = 5005000
$ws.Range("M2").Value = -5004888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 4000.5
$ws.Range("I20").Value = 4000
$ws.Range("K20").Value = 4000
$ws.Range("M20").Value = -3760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1386.7894
$ws.Range("I126").Value = 911.53845
$ws.Range("J126").Value = 2416.5
$ws.Range("K126").Value = 2734.61535
$ws.Range("L126").Value = 7249.5
$ws.Range("M126").Value = -264.61535
$ws.Range("N126").Value = -12189.5
